# Auto-generated edit script: updates leve profit calculation values
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1019.5
$ws.Range("J97").Value = 1036.875
$ws.Range("L97").Value = 3110.625
$ws.Range("N97").Value = -4102.625

$ws.Range("H112").Value = 4851.78
$ws.Range("J112").Value = 5112.5317
$ws.Range("L112").Value = 15337.5951
$ws.Range("N112").Value = -17553.5951

$ws.Range("H138").Value = 4455.5874
$ws.Range("I138").Value = 4393.4614
$ws.Range("J138").Value = 4471.74
$ws.Range("K138").Value = 13180.3842
$ws.Range("L138").Value = 13415.22
$ws.Range("M138").Value = -8040.3842
$ws.Range("N138").Value = -23695.22

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15167288
$ws.Range("I32").Value = 19621398
$ws.Range("J32").Value = 23318.934
$ws.Range("K32").Value = 19621398
$ws.Range("L32").Value = 23318.934
$ws.Range("M32").Value = -19621111
$ws.Range("N32").Value = -23892.934

$ws.Range("H74").Value = 10418650
$ws.Range("I74").Value = 1037.9395
$ws.Range("K74").Value = 1037.9395
$ws.Range("M74").Value = -163.9395

$ws.Range("H77").Value = 10418650
$ws.Range("I77").Value = 1037.9395
$ws.Range("K77").Value = 5189.6975
$ws.Range("M77").Value = -821.6975000000002

$ws.Range("H122").Value = 334504
$ws.Range("I122").Value = 500756
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 1502268
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1499818
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 1977433.4
$ws.Range("I132").Value = 5408.6294
$ws.Range("J132").Value = 6414489
$ws.Range("K132").Value = 16225.8882
$ws.Range("L132").Value = 19243467
$ws.Range("M132").Value = -13695.8882
$ws.Range("N132").Value = -19248527

$ws.Range("H133").Value = 49660
$ws.Range("J133").Value = 49660
$ws.Range("L133").Value = 49660
$ws.Range("N133").Value = -54720

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 57680
$ws.Range("J51").Value = 57680
$ws.Range("L51").Value = 57680
$ws.Range("N51").Value = -58662

$ws.Range("H88").Value = 35635.445
$ws.Range("J88").Value = 35635.445
$ws.Range("L88").Value = 35635.445
$ws.Range("N88").Value = -36447.445

$ws.Range("H91").Value = 35635.445
$ws.Range("J91").Value = 35635.445
$ws.Range("L91").Value = 35635.445
$ws.Range("N91").Value = -38443.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1331.1428
$ws.Range("I16").Value = 1414.5714
$ws.Range("J16").Value = 1164.2858
$ws.Range("K16").Value = 1414.5714
$ws.Range("L16").Value = 1164.2858
$ws.Range("M16").Value = -1127.5714
$ws.Range("N16").Value = -1738.2858

$ws.Range("H99").Value = 2441.15
$ws.Range("I99").Value = 1908
$ws.Range("J99").Value = 2535.2354
$ws.Range("K99").Value = 1908
$ws.Range("L99").Value = 2535.2354
$ws.Range("M99").Value = -410
$ws.Range("N99").Value = -5531.2354

$ws.Range("H113").Value = 1331.1428
$ws.Range("I113").Value = 1414.5714
$ws.Range("J113").Value = 1164.2858
$ws.Range("K113").Value = 1414.5714
$ws.Range("L113").Value = 1164.2858
$ws.Range("M113").Value = 755.4286
$ws.Range("N113").Value = -5504.2858

$ws.Range("H122").Value = 2407.125
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 2465.2856
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 7395.8568
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -12295.8568

$ws.Range("H126").Value = 2441.15
$ws.Range("I126").Value = 1908
$ws.Range("J126").Value = 2535.2354
$ws.Range("K126").Value = 5724
$ws.Range("L126").Value = 7605.706200000001
$ws.Range("M126").Value = -3254
$ws.Range("N126").Value = -12545.7062

$ws.Range("H132").Value = 19609724
$ws.Range("I132").Value = 23811016
$ws.Range("K132").Value = 71433048
$ws.Range("M132").Value = -71430518

$ws.Range("H141").Value = 74962.35000000001
$ws.Range("J141").Value = 74827.125
$ws.Range("L141").Value = 74827.125
$ws.Range("N141").Value = -85187.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 1814.2858
$ws.Range("I48").Value = 1566.6666
$ws.Range("K48").Value = 4699.9998
$ws.Range("M48").Value = -4449.9998

$ws.Range("H55").Value = 1992.0834
$ws.Range("J55").Value = 1992.0834
$ws.Range("L55").Value = 5976.2502
$ws.Range("N55").Value = -6330.2502

$ws.Range("H81").Value = 5728.5454
$ws.Range("I81").Value = 1571.4286
$ws.Range("J81").Value = 13003.5
$ws.Range("K81").Value = 4714.2858
$ws.Range("L81").Value = 39010.5
$ws.Range("M81").Value = -3591.2858
$ws.Range("N81").Value = -41256.5

$ws.Range("H84").Value = 5728.5454
$ws.Range("I84").Value = 1571.4286
$ws.Range("J84").Value = 13003.5
$ws.Range("K84").Value = 14142.8574
$ws.Range("L84").Value = 117031.5
$ws.Range("M84").Value = -8526.857399999999
$ws.Range("N84").Value = -128263.5

$ws.Range("H107").Value = 29413156
$ws.Range("J107").Value = 47621190
$ws.Range("L107").Value = 142863570
$ws.Range("N107").Value = -142867410

$ws.Range("H113").Value = 720.94116
$ws.Range("I113").Value = 710.8333
$ws.Range("J113").Value = 732.3125
$ws.Range("K113").Value = 2132.4999
$ws.Range("L113").Value = 2196.9375
$ws.Range("M113").Value = 37.5001000000002
$ws.Range("N113").Value = -6536.9375

$ws.Range("H131").Value = 3197
$ws.Range("I131").Value = 637.9167
$ws.Range("J131").Value = 3946
$ws.Range("K131").Value = 1913.7501
$ws.Range("L131").Value = 11838
$ws.Range("M131").Value = 3126.2499
$ws.Range("N131").Value = -21918

$ws.Range("H137").Value = 9293638
$ws.Range("J137").Value = 18577838
$ws.Range("L137").Value = 55733514
$ws.Range("N137").Value = -55743714

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1908.4348
$ws.Range("I102").Value = 1793.3334
$ws.Range("J102").Value = 2322.8
$ws.Range("K102").Value = 1793.3334
$ws.Range("L102").Value = 2322.8
$ws.Range("M102").Value = -171.3334
$ws.Range("N102").Value = -5566.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8035.857
$ws.Range("I22").Value = 875
$ws.Range("J22").Value = 10900.2
$ws.Range("K22").Value = 875
$ws.Range("L22").Value = 10900.2
$ws.Range("M22").Value = -580
$ws.Range("N22").Value = -11490.2

$ws.Range("H27").Value = 8035.857
$ws.Range("I27").Value = 875
$ws.Range("J27").Value = 10900.2
$ws.Range("K27").Value = 875
$ws.Range("L27").Value = 10900.2
$ws.Range("M27").Value = -768
$ws.Range("N27").Value = -11114.2

$ws.Range("H61").Value = 2042
$ws.Range("I61").Value = 888.1177
$ws.Range("K61").Value = 888.1177
$ws.Range("M61").Value = -686.1177

$ws.Range("H100").Value = 40146.04
$ws.Range("I100").Value = 85649.37
$ws.Range("J100").Value = 4393.4287
$ws.Range("K100").Value = 85649.37
$ws.Range("L100").Value = 4393.4287
$ws.Range("M100").Value = -85108.37
$ws.Range("N100").Value = -5475.4287

$ws.Range("H113").Value = 2042
$ws.Range("I113").Value = 888.1177
$ws.Range("K113").Value = 888.1177
$ws.Range("M113").Value = 1281.8823

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 480.5
$ws.Range("I107").Value = 480.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1441.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 478.5
$ws.Range("N107").ClearContents()

$ws.Range("H122").Value = 2847.3157
$ws.Range("I122").Value = 2406.6
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 7219.799999999999
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -4769.799999999999
$ws.Range("N122").Value = -18400
